$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Companies")

# Row 2: Walmart
$ws.Range("A2").Value = "Walmart"
$ws.Range("B2").Value = "Es una corporación multinacional de tiendas de origen estadounidense, que opera cadenas de grandes almacenes de descuento y clubes de almacenes."
$ws.Range("C2").Value = "Centro Comercial, 3 Avenida 41 10, Cdad. de Guatemala"
$ws.Range("D2").Value = "walmartgt@gmail.com"
$ws.Range("E2").Value = "HIGH_LEVEL"
$ws.Range("F2").Value = 62
$ws.Range("G2").Value = "Minorista comercial"
$ws.Range("H2").Value = "https://www.walmart.com.gt/?srsltid=AfmBOopQ2UETpvuX6b9aIQe3QiUt2oghiFLsLMzfdT8bWOd6YrKBDAlb"

# Row 3: Gallo
$ws.Range("A3").Value = "Gallo"
$ws.Range("B3").Value = "Cerveza Gallo es una marca de cerveza producida por la Cervecería Centro Americana, S.A.. Es la cerveza más consumida en Guatemala. La Cervecería Centro Americana tuvo hasta 2003 una cuota de mercado cercana al 100% del mercado cervecero guatemalteco."
$ws.Range("C3").Value = "3a Avenida Norte Final, Interior Finca El Zapote Z. 2, 3A Av · 1801 237 8392"
$ws.Range("D3").Value = "gallo@gmail.com"
$ws.Range("E3").Value = "HIGH_LEVEL"
$ws.Range("F3").Value = 99
$ws.Range("G3").Value = "Cerveceria"
$ws.Range("H3").Value = "https://www.elgallomasgallo.com.gt/?srsltid=AfmBOorr0gjphXxHUzznCHEjZBU5IjWEDzXufsf2nVPRmNhKFg91qWOF"

# Row 4: DollarCity
$ws.Range("A4").Value = "DollarCity"
$ws.Range("B4").Value = "Encuentra todo lo que necesitas en Decoración, Hogar, Oficina, Mascotas y mucho más. ¡Conoce nuestras ubicaciones entrando ahora!"
$ws.Range("C4").Value = "Calzada Roosevelt 9-12 · 2210 6000"
$ws.Range("D4").Value = "dollarcity@gmail.com"
$ws.Range("E4").Value = "MEDIUM_LEVEL"
$ws.Range("F4").Value = 32
$ws.Range("G4").Value = "Comercial de productos"
$ws.Range("H4").Value = "https://dollarcity.com/"

# Row 5: CocaCola
$ws.Range("A5").Value = "CocaCola"
$ws.Range("B5").Value = "The Coca-Cola Company es una corporación multinacional estadounidense de bebidas."
$ws.Range("C5").Value = "Bulevar Aguilar Batres · 2413 7500"
$ws.Range("D5").Value = "cocacola@gmail.com"
$ws.Range("E5").Value = "HIGH_LEVEL"
$ws.Range("F5").Value = 150
$ws.Range("G5").Value = "Alimentos y bebidas"
$ws.Range("H5").Value = "https://cocacola.com/"

# Row 6: PolloCampero
$ws.Range("A6").Value = "PolloCampero"
$ws.Range("B6").Value = "Vendemos pollo siiiii"
$ws.Range("C6").Value = "3ra Avenida, 5ta Calle, Alamierda el calvario"
$ws.Range("D6").Value = "pollocampero@gmail.com"
$ws.Range("E6").Value = "LOW_LEVEL"
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = "Alimentos y bebidas"
$ws.Range("H6").Value = "https://pollocampero.com/"
